# New Instances + Updated LPM
# Populates Sheet1 with the LPM instance-run summary table: a bold/boxed
# header row (A1:Q1) plus one data row (A2:Q2) for instance "vCp0n150s24d5".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$headers = @(
    "Instancia", "Memoria", "Status", "Valor FO", "Best Bound", "Rel GAP",
    "Tiempo Carga", "Tiempo Ejec", "Pacientes Atend", "Prioridad",
    "Avg Fichas", "Std Fichas", "Avg Cirug", "Std Cirug", "Avg Ratio",
    "Std Ratio", "Ocupación"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Style the header cells one at a time: bold font, thin box border, and
# centered / top-aligned text.
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Borders.LineStyle = 1
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
}

# ---- Data row ---------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "vCp0n150s24d5"
$ws.Cells.Item(2, 2).Value = 12500.55
$ws.Cells.Item(2, 3).Value = "Fact"
$ws.Cells.Item(2, 4).Value = 134746
$ws.Cells.Item(2, 5).Value = 138803.35
$ws.Cells.Item(2, 6).Value = 0.03

# "Tiempo Carga" is stored as text (looks numeric) rather than a number -
# a leading apostrophe forces Excel to keep it as a text value.
$ws.Cells.Item(2, 7).Value = "'517.99"

$ws.Cells.Item(2, 8).Value = 2389.33
$ws.Cells.Item(2, 9).Value = 143
$ws.Cells.Item(2, 10).Value = 143
$ws.Cells.Item(2, 11).Value = 78.82
$ws.Cells.Item(2, 12).Value = 41.49
$ws.Cells.Item(2, 13).Value = 6.5
$ws.Cells.Item(2, 14).Value = 3.3
$ws.Cells.Item(2, 15).Value = 12.18
$ws.Cells.Item(2, 16).Value = 2
$ws.Cells.Item(2, 17).Value = 0.56

# ---- Page setup (margins in points; 72pt = 1in) ------------------------
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72
